$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): shift years right by two columns and add 2025 ---
# Copy header style (bold + border + center/top alignment) onto the newly used columns G, H, I
$ws.Range("F1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B1").Value = 2018
$ws.Range("C1").Value = 2019
$ws.Range("D1").Value = 2020
$ws.Range("E1").Value = 2021
$ws.Range("F1").Value = 2022
$ws.Range("G1").Value = 2023
$ws.Range("H1").Value = 2024
$ws.Range("I1").Value = 2025

# --- Data rows 2-13 (months 1-12): refreshed figures, 2018/2019/2025 columns added or removed per row ---
# Row 2 (month 1)
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 70181.04
$ws.Range("D2").Value = 273853.03
$ws.Range("E2").Value = 505264.48
$ws.Range("F2").Value = 444724.53
$ws.Range("G2").Value = 528679.07
$ws.Range("H2").Value = 454942.61
$ws.Range("I2").Value = 1077628.61

# Row 3 (month 2)
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 55926.21
$ws.Range("D3").Value = 224777.98
$ws.Range("E3").Value = 308589.64
$ws.Range("F3").Value = 423121.33
$ws.Range("G3").Value = 557000.89
$ws.Range("H3").Value = 486371.27
$ws.Range("I3").Value = 941184.13

# Row 4 (month 3)
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value = 58748
$ws.Range("D4").Value = 230762.37
$ws.Range("E4").Value = 471550.93
$ws.Range("F4").Value = 461440.4
$ws.Range("G4").Value = 661369.18
$ws.Range("H4").Value = 570778.23
$ws.Range("I4").Value = 796983.15

# Row 5 (month 4)
$ws.Range("B5").ClearContents()
$ws.Range("C5").Value = 69525
$ws.Range("D5").Value = 231723.38
$ws.Range("E5").Value = 414500.83
$ws.Range("F5").Value = 528483.57
$ws.Range("G5").Value = 592203.52
$ws.Range("H5").Value = 553985.89
$ws.Range("I5").Value = 449688.27

# Row 6 (month 5)
$ws.Range("B6").ClearContents()
$ws.Range("C6").Value = 58295.7
$ws.Range("D6").Value = 298073.29
$ws.Range("E6").Value = 457351.79
$ws.Range("F6").Value = 597756.95
$ws.Range("G6").Value = 684242.49
$ws.Range("H6").Value = 855536.36
$ws.Range("I6").ClearContents()

# Row 7 (month 6)
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = 75018.6
$ws.Range("D7").Value = 283819.25
$ws.Range("E7").Value = 504810.11
$ws.Range("F7").Value = 520962.54
$ws.Range("G7").Value = 497416.03
$ws.Range("H7").Value = 625142.05
$ws.Range("I7").ClearContents()

# Row 8 (month 7)
$ws.Range("B8").Value = 75489.75
$ws.Range("C8").Value = 102664.2
$ws.Range("D8").Value = 291646.71
$ws.Range("E8").Value = 814721.26
$ws.Range("F8").Value = 607227.32
$ws.Range("G8").Value = 614574.44
$ws.Range("H8").Value = 950349.51
$ws.Range("I8").ClearContents()

# Row 9 (month 8)
$ws.Range("B9").Value = 180862
$ws.Range("C9").Value = 195298.1
$ws.Range("D9").Value = 322356.9
$ws.Range("E9").Value = 627045.92
$ws.Range("F9").Value = 551353.81
$ws.Range("G9").Value = 510615.64
$ws.Range("H9").Value = 850701.4
$ws.Range("I9").ClearContents()

# Row 10 (month 9)
$ws.Range("B10").Value = 107828.15
$ws.Range("C10").Value = 175841.69
$ws.Range("D10").Value = 285208.92
$ws.Range("E10").Value = 620390.4
$ws.Range("F10").Value = 617174.95
$ws.Range("G10").Value = 510107.45
$ws.Range("H10").Value = 716637.28
$ws.Range("I10").ClearContents()

# Row 11 (month 10)
$ws.Range("B11").Value = 99858.93
$ws.Range("C11").Value = 299112.23
$ws.Range("D11").Value = 268079.69
$ws.Range("E11").Value = 577095.91
$ws.Range("F11").Value = 828838.14
$ws.Range("G11").Value = 612264.49
$ws.Range("H11").Value = 1064675.22
$ws.Range("I11").ClearContents()

# Row 12 (month 11)
$ws.Range("B12").Value = 105543.5
$ws.Range("C12").Value = 259819.03
$ws.Range("D12").Value = 272246.01
$ws.Range("E12").Value = 731062.09
$ws.Range("F12").Value = 736765.16
$ws.Range("G12").Value = 490239.54
$ws.Range("H12").Value = 1031789.3
$ws.Range("I12").ClearContents()

# Row 13 (month 12)
$ws.Range("B13").Value = 130615
$ws.Range("C13").Value = 321287.28
$ws.Range("D13").Value = 331309.5
$ws.Range("E13").Value = 532112.76
$ws.Range("F13").Value = 1042718.34
$ws.Range("G13").Value = 674121.83
$ws.Range("H13").Value = 1264588.39
$ws.Range("I13").ClearContents()
